# Final fix of weibull unit case for testing.
#
# Applies the data + view-state changes described by the commit:
#  - Capacity_start!D2, D3: 5 -> 20
#  - Capacity_new!F2,F6,...,F238 (every 4th row): 0 -> 2
#  - Tech_lifetime!E2:E5: 15/10/20/15 -> 10/5/15/10
#  - View state: Tech_lifetime's selection moves to G10 and is no longer
#    the active tab; Capacity_new becomes the active tab, its view
#    scrolls/selects L23.

$wb = $excel.ActiveWorkbook

# --- Capacity_start sheet: D2 and D3 5 -> 20 ---
$wsCapacityStart = $wb.Worksheets.Item("Capacity_start")
$wsCapacityStart.Range("D2").Value = 20
$wsCapacityStart.Range("D3").Value = 20

# --- Capacity_new sheet: column F (row 2, 6, 10, ... 238) 0 -> 2 ---
$wsCapacityNew = $wb.Worksheets.Item("Capacity_new")
for ($r = 2; $r -le 238; $r += 4) {
    $wsCapacityNew.Cells.Item($r, 6).Value = 2
}

# --- Tech_lifetime sheet: E2:E5 values updated ---
$wsTechLifetime = $wb.Worksheets.Item("Tech_lifetime")
$wsTechLifetime.Range("E2").Value = 10
$wsTechLifetime.Range("E3").Value = 5
$wsTechLifetime.Range("E4").Value = 15
$wsTechLifetime.Range("E5").Value = 10

# Update Tech_lifetime's selection first (while it's still the active
# sheet) so it ends up non-active once Capacity_new is activated below.
$wsTechLifetime.Range("G10").Select()

# --- Capacity_new becomes the active sheet/tab with a new selection ---
$wsCapacityNew.Activate()
$wsCapacityNew.Range("L23").Select()
